$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Permendagri No.24 Tahun 2021" -> "Permendagri Nomor 47 Tahun 2021"
#    The phrase is split across two runs in the source paragraph:
#      run A: "Permendagri No.24 Tahun 202"   (keeps its w:rsidRPr)
#      run B: "1"                             (gets removed entirely)
#    After the edit there must be a single run carrying the full new text
#    and the original formatting / rsid of run A.
# ---------------------------------------------------------------------------

$old = "Permendagri No.24 Tahun 2021"
$new = "Permendagri Nomor 47 Tahun 2021"

$locate = $d.Content
$found = $locate.Find.Execute($old, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $locate.Start
    $matchEnd   = $locate.End

    # The trailing "1" lives in its own run that is exactly the last
    # character of the match; drop it first so the preceding run keeps its
    # original <w:r> attributes (e.g. w:rsidRPr) untouched.
    $tailRun = $d.Range($matchEnd - 1, $matchEnd)
    $tailRun.Delete()

    # Re-insert the remainder of the paragraph's first run as a tiny OOXML
    # fragment so its rPr / rsid attributes are preserved exactly, only the
    # visible text changes.
    $headRun = $d.Range($matchStart, $matchEnd - 1)

    $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
               '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
               '<w:r w:rsidRPr="009E259F"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:lang w:val="en-US"/></w:rPr><w:t>' + $new + '</w:t></w:r>' + `
               '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $headRun.InsertXML($xmlFrag)
}
